$wb = $excel.ActiveWorkbook

$wsReg = $wb.Worksheets.Item(1)
[void]$wsReg.Range("E4").Select()

# --- Add "Variables" sheet after RegisterNodes ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVars = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsVars.Name = "Variables"

$wsVars.Range("B1").Value = "Tenant Name"
$wsVars.Range("A1").Value = "Apic Address"
$wsVars.Range("C1").Value = "PodId"
$wsVars.Range("A2").Value = "sandboxapicdc.cisco.com"
$wsVars.Range("B2").Value = "TestTenantName"
$wsVars.Range("C2").Value = 1
$wsVars.Range("C3").Value = 2

[void]$wsVars.Range("E3").Select()

# --- Add "test" sheet after Variables ---
$last2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTest = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last2)
$wsTest.Name = "test"

$wsTest.Range("A1").Value = "Apic"
$wsTest.Range("B1").Value = "Tenant"
$wsTest.Range("C1").Value = "PodId"

# --- Data validations on "test" sheet, sourced from Variables ---
[void]$wsTest.Range("C3:C14").Validation.Add(3, 1, 1, '=Variables!$C$2:$C$3')
[void]$wsTest.Range("C2").Validation.Add(3, 1, 1, '=Variables!$C$2:$C$22')
[void]$wsTest.Range("B2").Validation.Add(3, 1, 1, '=Variables!$B$2:$B$60')
[void]$wsTest.Range("A2").Validation.Add(3, 1, 1, '=Variables!$A$2:$A$50')

[void]$wsTest.Range("I9").Select()

Write-Host "Worksheets:" $wb.Worksheets.Count
